$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new values look like plain numbers need to be forced to
# Text format first, so Excel keeps them as strings (matching the source data which
# stores prices/percentages as text, including values with trailing zeros).
$textCells = @("D5", "D6", "D8", "D11", "D14", "D18", "D20", "D21", "D23", "D24", "D27", "D31", "D32", "D35", "D36", "D38", "D39", "D42", "D43", "D44", "D46", "D48", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "63.374.00"
$ws.Range("E2").Value = "  -1.21%  "
$ws.Range("D3").Value = "2.712.93"
$ws.Range("E3").Value = "  -1.74%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "560.48"
$ws.Range("E5").Value = "  -2.98%  "
$ws.Range("D6").Value = "156.53"
$ws.Range("E6").Value = "  -1.76%  "
$ws.Range("D8").Value = "0.589"
$ws.Range("E8").Value = "  -2.56%  "
$ws.Range("E9").Value = "  -3.02%  "
$ws.Range("E10").Value = "  +0.57%  "
$ws.Range("D11").Value = "5.59"
$ws.Range("E11").Value = "  -1.62%  "
$ws.Range("E12").Value = "  -4.38%  "
$ws.Range("D13").Value = "3.193.87"
$ws.Range("E13").Value = "  -1.69%  "
$ws.Range("D14").Value = "26.35"
$ws.Range("E14").Value = "  -2.26%  "
$ws.Range("D15").Value = "63.226.80"
$ws.Range("E15").Value = "  -0.80%  "
$ws.Range("E16").Value = "  -3.37%  "
$ws.Range("D17").Value = "2.713.95"
$ws.Range("E17").Value = "  -1.82%  "
$ws.Range("D18").Value = "12.16"
$ws.Range("E18").Value = "  +0.08%  "
$ws.Range("E19").Value = "  -4.33%  "
$ws.Range("D20").Value = "350.92"
$ws.Range("E20").Value = "  -1.85%  "
$ws.Range("D21").Value = "6.43"
$ws.Range("E21").Value = "  -4.82%  "
$ws.Range("E22").Value = "  -0.27%  "
$ws.Range("D23").Value = "0.510"
$ws.Range("E23").Value = "  -4.66%  "
$ws.Range("D24").Value = "64.16"
$ws.Range("E24").Value = "  -2.12%  "
$ws.Range("E25").Value = "  -1.22%  "
$ws.Range("E26").Value = "  +0.17%  "
$ws.Range("D27").Value = "8.19"
$ws.Range("E27").Value = "  -4.93%  "
$ws.Range("E28").Value = "  -2.51%  "
$ws.Range("E29").Value = "  +9.18%  "
$ws.Range("E30").Value = "  -0.74%  "
$ws.Range("D31").Value = "7.16"
$ws.Range("E31").Value = "  -2.24%  "
$ws.Range("D32").Value = "165.73"
$ws.Range("E32").Value = "  -2.30%  "
$ws.Range("E33").Value = "  -1.26%  "
$ws.Range("E34").Value = "  +0.03%  "
$ws.Range("D35").Value = "19.80"
$ws.Range("E35").Value = "  -2.41%  "
$ws.Range("D36").Value = "4.81"
$ws.Range("E36").Value = "  -2.53%  "
$ws.Range("E37").Value = "  -2.87%  "
$ws.Range("D38").Value = "344.48"
$ws.Range("E38").Value = "  -0.33%  "
$ws.Range("D39").Value = "0.960"
$ws.Range("E39").Value = "  -4.84%  "
$ws.Range("E40").Value = "  -3.82%  "
$ws.Range("E41").Value = "  -4.27%  "
$ws.Range("D42").Value = "38.47"
$ws.Range("D43").Value = "21.37"
$ws.Range("E43").Value = "  -2.32%  "
$ws.Range("D44").Value = "20.72"
$ws.Range("E44").Value = "  -3.44%  "
$ws.Range("E45").Value = "  -3.42%  "
$ws.Range("D46").Value = "0.623"
$ws.Range("E46").Value = "  -1.46%  "
$ws.Range("E47").Value = "  +0.12%  "
$ws.Range("D48").Value = "131.87"
$ws.Range("E48").Value = "  -2.91%  "
$ws.Range("E49").Value = "  +0.09%  "
$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D50").Value = "0.0245"
$ws.Range("E50").Value = "  -4.13%  "
$ws.Range("B51").Value = "Stellar"
$ws.Range("C51").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D51").Value = "0.0982"
$ws.Range("E51").Value = "  -3.92%  "
